$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SumProduct")

# Add new cells for the "Ignored cells" test scenario.
# Cells are written in the order that reproduces the shared-string table
# append order recorded in the workbook (Ignored, Will, be, Ignored after
# two empty cells).
$ws.Range("H7").Value = "Ignored"
$ws.Range("F5").Value = "Will"
$ws.Range("G6").Value = "be"
$ws.Range("F4").Value = "Ignored after two empty cells"

# Set column C width to match the bestFit width seen in sheet1 (15.28515625)
$ws.Columns.Item(3).EntireColumn.AutoFit()

# Update the selection to match the recorded state after editing.
$ws.Range("G15").Select() | Out-Null
